$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")

# Rename the contribution-assessment-limit rows: the old split between
# "GKV Ost" / "GKV West" is replaced by a single "Beitragsbemessungsgrenze GKV"
# row, and a new "Jahresarbeitsentgeltgrenze GKV" row.
$ws1.Range("A5").Value = "Beitragsbemessungsgrenze GKV"
$ws1.Range("A6").Value = "Jahresarbeitsentgeltgrenze GKV"

# Move/Update the active selection on Tabelle1 to A9, as recorded in the
# saved workbook view.
$ws1.Activate()
$ws1.Range("A9").Select()
